$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '37.864.56'
$c.ClearFormats()
$ws.Range('E2').Value = '  +0.08%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.030.97'
$c.ClearFormats()
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  -0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '227.43'
$c.ClearFormats()
$ws.Range('E5').Value = '  -0.91%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.615'
$c.ClearFormats()
$ws.Range('E6').Value = '  -0.17%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '59.61'
$c.ClearFormats()
$ws.Range('E7').Value = '  +2.41%  '
$ws.Range('E8').Value = '  -0.02%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.383'
$c.ClearFormats()
$ws.Range('E9').Value = '  -0.58%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.0810'
$c.ClearFormats()
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('E12').Value = '  +0.05%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '2.330.40'
$c.ClearFormats()
$ws.Range('E13').Value = '  -0.80%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '21.15'
$c.ClearFormats()
$ws.Range('E14').Value = '  +2.28%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.761'
$c.ClearFormats()
$ws.Range('E15').Value = '  +1.67%  '
$ws.Range('E16').Value = '  -1.68%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '2.049.13'
$c.ClearFormats()
$ws.Range('E17').Value = '  -0.27%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '37.768.04'
$c.ClearFormats()
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('E19').Value = '  -2.14%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '70.04'
$c.ClearFormats()
$ws.Range('E20').Value = '  +0.57%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '0.0₃0824'
$c.ClearFormats()
$ws.Range('E21').Value = '  -0.91%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '224.99'
$c.ClearFormats()
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -2.19%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.ClearFormats()
$ws.Range('E25').Value = '  -1.70%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '9.28'
$c.ClearFormats()
$ws.Range('E26').Value = '  +0.24%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '165.34'
$c.ClearFormats()
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('E28').Value = '  -2.59%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '18.96'
$c.ClearFormats()
$ws.Range('E29').Value = '  -0.32%  '
$ws.Range('E30').Value = '  -4.37%  '
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('B32').Value = 'WEMIXToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '2.13'
$c.ClearFormats()
$ws.Range('E32').Value = '  +2.05%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '4.44'
$c.ClearFormats()
$ws.Range('E33').Value = '  -2.05%  '
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('E35').Value = '  -1.44%  '
$ws.Range('E36').Value = '  +6.61%  '
$ws.Range('E37').Value = '  -3.06%  '
$ws.Range('E38').Value = '  -2.31%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E39').Value = '  -0.12%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '1.520.61'
$c.ClearFormats()
$ws.Range('E40').Value = '  +2.41%  '
$ws.Range('E41').Value = '  +0.72%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '96.68'
$c.ClearFormats()
$ws.Range('E42').Value = '  -0.76%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '16.80'
$c.ClearFormats()
$ws.Range('E43').Value = '  +1.34%  '
$ws.Range('E44').Value = '  -0.35%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.0916'
$c.ClearFormats()
$ws.Range('E45').Value = '  -1.91%  '
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('E47').Value = '  -3.44%  '
$ws.Range('E48').Value = '  -0.48%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.ClearFormats()
$ws.Range('E49').Value = '  +0.14%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '7.09'
$c.ClearFormats()
$ws.Range('E50').Value = '  +1.37%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '2.218.49'
$c.ClearFormats()
$ws.Range('E51').Value = '  -0.98%  '
